# ---------------------------------------------------------------------------
# Edit: SRS/2.4.1.docx
# Replaces the first body paragraph's (P2) content with the new intro text,
# inserts a brand-new paragraph (P3) describing the two-module pipeline
# (carrying the relocated "_GoBack" bookmark at its end), empties the old
# "因此，本系統將按..." paragraph (keeping only its paragraph mark / pPr,
# now with an eastAsia hint on its rFonts), and removes the figure paragraph,
# the figure-caption paragraph and the trailing bookmark-only paragraph.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. Replace paragraph 2 (old intro) with the new P2 + new P3 runs ------
$p2 = $d.Paragraphs.Item(2)
$insertAt = $p2.Range
$insertAt.Collapse(1)
$p2.Range.Delete()

$newP2P3 = '<w:p ' + $wNs + '><w:pPr><w:ind w:firstLine="480"/><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>本系統的主要目標，是提出一套</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體"/></w:rPr><w:t>基於深度網路之人臉情感檢測系統</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>，將拍攝後的從業人員影像，運用人工智慧(Artificial Intelligence，AI)中的深度學習(Deep Learning)，配合卷積神經網路(Convolution Neural Networks，CNNs)進行情感檢測，提供服務產業在服務態度的自動化管理系統。</w:t></w:r></w:p><w:p ' + $wNs + '><w:pPr><w:ind w:firstLine="480"/><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>本系統</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>分為兩部分</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>實施</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>，</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>第一部分是人臉偵測模組，</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>第二</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>部分</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>是情緒量化模組</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>。</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>首先，將</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>攝影機</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>或監視器</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>拍攝的圖片</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>，</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>輸入至人臉偵測模組，獲得該從業人員的臉部位置，並將它裁切。之後將裁切的圖片輸入至情緒量化模組，進行笑容量化數值的計算</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>，並視覺化</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>顯示於監控畫面，</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微軟正黑體" w:eastAsia="微軟正黑體" w:hAnsi="微軟正黑體" w:hint="eastAsia"/></w:rPr><w:t>具體呈現從業人員在服務態度上的情緒表現。</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$insertAt.InsertXML($newP2P3)

# --- 2. Empty the old "因此，本系統將按..." paragraph (now #4), but keep  --
#        its paragraph mark, switching its pPr rFonts to carry the eastAsia
#        hint as in the target revision.
$p4 = $d.Paragraphs.Item(4)
$p4Insert = $p4.Range
$p4Insert.Collapse(1)
$p4.Range.Delete()

$newP4 = '<w:p ' + $wNs + '><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'
$p4Insert.InsertXML($newP4)

# --- 3. Remove the figure paragraph (image), the figure-caption paragraph -
#        and the trailing bookmark-only paragraph (now #5, #5, #5 since
#        each deletion shifts the following ones up).
$d.Paragraphs.Item(5).Range.Delete()
$d.Paragraphs.Item(5).Range.Delete()
$d.Paragraphs.Item(5).Range.Delete()

Write-Output ("FinalParagraphCount=" + $d.Paragraphs.Count)
